# Add a new daily check-in row for 2022-05-25 (LeetCode "替换空格") into row 6
# of the "打卡模板" sheet, matching the pattern of rows 4/5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("打卡模板")

# Row 6 columns: A=Date, B=Title, C=Difficulty, D=Status, E=Time(min),
# F=Category, G=Summary, H=NeedRedo?, I=NewKnowledge, J=Status2, K=Redo?
# The date column stores plain text (e.g. "2022.5.23"), not a real date
# serial, so force text formatting before assigning it.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2022.5.25"
$ws.Range("B6").Value = "剑指offer 替换空格（https://leetcode.cn/problems/ti-huan-kong-ge-lcof/）"
$ws.Range("C6").Value = "简单"
$ws.Range("D6").Value = "通过"
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = "字符串"
$ws.Range("G6").Value = "1.字符串初始化等。"
$ws.Range("H6").Value = "待巩固"
$ws.Range("I6").Value = "string的使用"
$ws.Range("J6").Value = "待学习"
$ws.Range("K6").Value = "是"

# Match the style used by the neighboring rows (A4/A5 -> A6 already has it).
$ws.Range("A6").Style = $ws.Range("A5").Style

# Update the view scroll position / selection to reflect the author's
# final cursor location after entering the new row.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("F30").Select()
